# Dev count run 20251205
$wb = $excel.ActiveWorkbook

# Update the summary sheet values/time stamp.
$ws = $wb.Worksheets.Item("SCM Report Summary")
$ws.Range("B3").Value = "10:16:38 AM"
$ws.Range("B5").Value = 0
$ws.Range("B7").Value = 0

# Remove the GitHub Details sheets entirely - no longer part of the report.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("GitHub Details").Delete()
$wb.Worksheets.Item("GitHub Details - Removed").Delete()
$excel.DisplayAlerts = $true
